# A new weekly price record ("2023-03-07") is inserted as row 24 of the
# "Hortaliza, Femacal de La Calera - Perejil" sheet; every existing record
# previously on rows 24-33 shifts down by one row (to rows 25-34), retaining
# its own values and formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data (old rows 24-33) down one row, carrying formatting
# (e.g. the date-formatted style on column D) along with it.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new observation.
$ws.Range("A24").Value = 3
$ws.Range("B24").Value = "Femacal de La Calera"
$ws.Range("C24").Value = "Coquimbo"
$ws.Range("D24").Value = 44992
$ws.Range("E24").Value = 5
$ws.Range("F24").Value = 100112044
$ws.Range("G24").Value = "Perejil"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 45
$ws.Range("K24").Value = 4000
$ws.Range("L24").Value = 4000
$ws.Range("M24").Value = 4000
$ws.Range("N24").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O24").Value = "Provincia de Quillota"
$ws.Range("P24").Value = 1333
$ws.Range("Q24").Value = 3
$ws.Range("R24").Value = "Hortaliza"
